$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows (20:21) pushing the signature block (old rows 24:25) down to 26:27
$ws.Range("B20:J21").Insert()

# Copy formatting from the existing "2508" period rows (18:19) into the new rows
$ws.Range("B18:J19").Copy()
$ws.Range("B20:J21").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B18:J19").Copy()
$ws.Range("B20:J21").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = 0

# Update the period for the newly added rows to "2509"
$ws.Range("E20").Value = "2509"
$ws.Range("E21").Value = "2509"

# Update the totals: one more period now included
$ws.Range("F13").Value = 3
$ws.Range("E11").Value = 343380
